# ---------------------------------------------------------------------------
# Applies the Dec-2023 data refresh described by the commit "Add files via
# upload" to data202312.xlsx:
#   * Sheet "部门情况202312"      -> updated F/G/H figures for 3 rows + a new
#                                    "普惠业务二部" row (row 10)
#   * Sheet "对公业务台账202312"  -> new "户均" column (U)
#   * Sheet "对公产品台账202312"  -> new "户均" / "平均利率" columns (J, K)
#   * Sheet "个人经营贷202312"    -> two rows removed (创业经营贷, 车商贷),
#                                    three new columns added (个人户均,
#                                    逾期金额, 不良金额), totals recomputed
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Helper: write a value that *looks* numeric ("12.00") but that the source
# workbook stores as literal text (inlineStr/shared-string) -- force the "@"
# text format before the assignment so the engine doesn't coerce it to a
# real number, then drop back to the Normal style so no stray number format
# sticks around on the cell.
function Set-TextValue($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# ===========================================================================
# Sheet 1: 部门情况202312
# ===========================================================================
$ws1 = $wb.Worksheets.Item("部门情况202312")

Set-TextValue $ws1 "F3" "2980.00"
Set-TextValue $ws1 "G3" "12.00"
Set-TextValue $ws1 "H3" "5.90"

Set-TextValue $ws1 "F4" "2409.71"
Set-TextValue $ws1 "G4" "45.00"
Set-TextValue $ws1 "H4" "6.04"

Set-TextValue $ws1 "F7" "251895.81"
Set-TextValue $ws1 "G7" "2885.00"
Set-TextValue $ws1 "H7" "6.09"

Set-TextValue $ws1 "F9" "92340.67"
Set-TextValue $ws1 "G9" "5691.00"
Set-TextValue $ws1 "H9" "17.64"

# New row 10: 普惠业务二部, all-zero figures (copy row 9's formatting down
# first so the new row looks like the rest of the table).
$ws1.Range("A9:O9").Copy()
$ws1.Range("A10:O10").PasteSpecial(-4122)

Set-TextValue $ws1 "A10" "普惠业务二部"
$cols1 = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O")
foreach ($c in $cols1) {
    $addr1 = $c + "10"
    Set-TextValue $ws1 $addr1 "0.00"
}

# ===========================================================================
# Sheet 3: 对公业务台账202312  -- add "户均" column (U)
# ===========================================================================
$ws3 = $wb.Worksheets.Item("对公业务台账202312")

# Match formatting of the existing last header/data column (T) before
# writing values.
$ws3.Range("T1").Copy()
$ws3.Range("U1").PasteSpecial(-4122)
$ws3.Range("U1").Value = "户均"

$u3 = @{
    "U2"  = "7776.91"
    "U3"  = "10786.67"
    "U4"  = "691.19"
    "U5"  = "218.44"
    "U6"  = "97.82"
    "U7"  = "569.41"
    "U8"  = "29.31"
    "U9"  = "595.09"
    "U10" = "1.82"
    "U11" = "33000.00"
    "U12" = "53766.66"
}
foreach ($addr in $u3.Keys) {
    Set-TextValue $ws3 $addr $u3[$addr]
}

# ===========================================================================
# Sheet 4: 对公产品台账202312  -- add "户均" (J) / "平均利率" (K) columns
# ===========================================================================
$ws4 = $wb.Worksheets.Item("对公产品台账202312")

$ws4.Range("I1").Copy()
$ws4.Range("J1:K1").PasteSpecial(-4122)
$ws4.Range("J1").Value = "户均"
$ws4.Range("K1").Value = "平均利率"

$j4 = @{
    2  = 0
    3  = 10786.67
    4  = 691.1900000000001
    5  = 218.44
    6  = 97.81999999999999
    7  = 300.46
    8  = 601.6
    9  = 29.31
    10 = 791.67
    11 = 160.76
    12 = 1000
    13 = 909.09
    14 = 509.56
    15 = 377.04
    16 = 1.82
    17 = 33000
    18 = 49475.43
}
$k4 = @{
    2  = 0
    3  = 5.51
    4  = 6.6
    5  = 6.6
    6  = 5.4
    7  = 6.8
    8  = 6.51
    9  = 6.3
    10 = 3.95
    11 = 5
    12 = 5.15
    13 = 5.15
    14 = 5
    15 = 6.5
    16 = 24
    17 = 7.5
    18 = 105.97
}
foreach ($r in $j4.Keys) {
    $ws4.Range("J$r").Value = $j4[$r]
    $ws4.Range("K$r").Value = $k4[$r]
}

# ===========================================================================
# Sheet 5: 个人经营贷202312
#   - drop row 2 (创业经营贷) and (the now-shifted) row for 车商贷
#   - add columns K "个人户均", L "逾期金额", M "不良金额"
#   - Total row recomputed to exclude the dropped rows
# ===========================================================================
$ws5 = $wb.Worksheets.Item("个人经营贷202312")

# 创业经营贷 is the first data row (row 2) -- remove it entirely.
$ws5.Rows.Item(2).Delete()
# 车商贷 was row 16 before the deletion above; after the shift it is row 15,
# directly above the Total row -- remove it too.
$ws5.Rows.Item(15).Delete()

# New headers (match formatting of the existing last header column, J).
$ws5.Range("J1").Copy()
$ws5.Range("K1:M1").PasteSpecial(-4122)
$ws5.Range("K1").Value = "个人户均"
$ws5.Range("L1").Value = "逾期金额"
$ws5.Range("M1").Value = "不良金额"

# New column values per remaining product row (2-14) plus the recomputed
# Total row (15).
$k5 = @{
    2  = 350
    3  = 245.32
    4  = 20.46
    5  = 37.9
    6  = 0.73
    7  = 4.87
    8  = 54.06
    9  = 68.59999999999999
    10 = 15.59
    11 = 50.18
    12 = 13.48
    13 = 66.48
    14 = 17.21
    15 = 944.88
}
$l5 = @{
    2  = 0
    3  = 1525
    4  = 0
    5  = 1736.68
    6  = 1.54
    7  = 0
    8  = 0
    9  = 0
    10 = 3809.28
    11 = 0
    12 = 8.52
    13 = 0
    14 = 17.21
    15 = 7098.23
}
$m5 = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 1059.38
    11 = 0
    12 = 0
    13 = 0
    14 = 17.21
    15 = 1076.59
}
foreach ($r in $k5.Keys) {
    $ws5.Range("K$r").Value = $k5[$r]
    $ws5.Range("L$r").Value = $l5[$r]
    $ws5.Range("M$r").Value = $m5[$r]
}

# Recompute the Total row (row 15) for the columns impacted by the two
# deleted rows.
$ws5.Range("B15").Value = 29172
$ws5.Range("C15").Value = 375372.3
$ws5.Range("D15").Value = 34069
$ws5.Range("E15").Value = 426924.73
$ws5.Range("F15").Value = 118.52
